$d = $word.ActiveDocument

# --- Locate the paragraph that holds the lone leading-space run right
# --- before the bold "Client vs. Server" run, and replace that space
# --- run with a (relocated) "_GoBack" bookmark.
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Client vs. Server*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $r = $target.Range
    $spaceRange = $d.Range($r.Start, $r.Start + 1)
    if ($spaceRange.Text -eq " ") {
        $spacePos = $spaceRange.Start
        $spaceRange.Text = ""
        $bmRange = $d.Range($spacePos, $spacePos)
        # Adding a bookmark with a name that already exists moves it,
        # which also removes the old "_GoBack" bookmark near the TCP
        # paragraph further down in the document.
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}
